$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a tiny precision drift in the previous day's timestamp (A12)
$ws.Range("A12").Value = 45812.39352998843

# Append the new day's price row
$ws.Range("A13").Value = 45813.39355503808
$ws.Range("B13").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C13").Value = "1Kg"
$ws.Range("D13").Value = "15,41€"

# Match the date-time number format used by the rest of column A
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
